$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1973.6666
$ws.Range("I32").Value = 2147.5
$ws.Range("J32").Value = 1886.75
$ws.Range("K32").Value = 2147.5
$ws.Range("L32").Value = 1886.75
$ws.Range("M32").Value = -1821.5
$ws.Range("N32").Value = -2538.75
$ws.Range("H53").Value = 961.1539
$ws.Range("I53").Value = 1950.1666
$ws.Range("J53").Value = 113.42857
$ws.Range("K53").Value = 1950.1666
$ws.Range("L53").Value = 113.42857
$ws.Range("M53").Value = -1313.1666
$ws.Range("N53").Value = -1387.42857
$ws.Range("H76").Value = 5099.7
$ws.Range("I76").Value = 3799.4
$ws.Range("K76").Value = 3799.4
$ws.Range("M76").Value = -3484.4
$ws.Range("H79").Value = 5099.7
$ws.Range("I79").Value = 3799.4
$ws.Range("K79").Value = 3799.4
$ws.Range("M79").Value = -2707.4
$ws.Range("H116").Value = 3031.7666
$ws.Range("I116").Value = 2604.2727
$ws.Range("J116").Value = 3279.2632
$ws.Range("K116").Value = 2604.2727
$ws.Range("L116").Value = 3279.2632
$ws.Range("M116").Value = 837.7273
$ws.Range("N116").Value = -10163.2632
$ws.Range("H125").Value = 5250
$ws.Range("I125").Value = 5250
$ws.Range("K125").Value = 47250
$ws.Range("M125").Value = -44790
$ws.Range("H137").Value = 13334296
$ws.Range("I137").Value = 828.75
$ws.Range("J137").Value = 28572544
$ws.Range("K137").Value = 2486.25
$ws.Range("L137").Value = 85717632
$ws.Range("M137").Value = 63.75
$ws.Range("N137").Value = -85722732
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7693.602
$ws.Range("I32").Value = 7705.3555
$ws.Range("J32").Value = 7641.0586
$ws.Range("K32").Value = 7705.3555
$ws.Range("L32").Value = 7641.0586
$ws.Range("M32").Value = -7418.3555
$ws.Range("N32").Value = -8215.0586
$ws.Range("H80").Value = 21086
$ws.Range("J80").Value = 21086
$ws.Range("L80").Value = 21086
$ws.Range("N80").Value = -23082
$ws.Range("H83").Value = 21086
$ws.Range("J83").Value = 21086
$ws.Range("L83").Value = 63258
$ws.Range("N83").Value = -73242
$ws.Range("H122").Value = 9538.357
$ws.Range("I122").Value = 14978
$ws.Range("J122").Value = 2285.5
$ws.Range("K122").Value = 44934
$ws.Range("L122").Value = 6856.5
$ws.Range("M122").Value = -42484
$ws.Range("N122").Value = -11756.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14609.3
$ws.Range("I82").Value = 11428.333
$ws.Range("J82").Value = 19380.75
$ws.Range("K82").Value = 11428.333
$ws.Range("L82").Value = 19380.75
$ws.Range("M82").Value = -11045.333
$ws.Range("N82").Value = -20146.75
$ws.Range("H85").Value = 14609.3
$ws.Range("I85").Value = 11428.333
$ws.Range("J85").Value = 19380.75
$ws.Range("K85").Value = 11428.333
$ws.Range("L85").Value = 19380.75
$ws.Range("M85").Value = -10102.333
$ws.Range("N85").Value = -22032.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6805048.5
$ws.Range("I31").Value = 1928.5834
$ws.Range("J31").Value = 25644456
$ws.Range("K31").Value = 1928.5834
$ws.Range("L31").Value = 25644456
$ws.Range("M31").Value = -1633.5834
$ws.Range("N31").Value = -25645046
$ws.Range("H34").Value = 6805048.5
$ws.Range("I34").Value = 1928.5834
$ws.Range("J34").Value = 25644456
$ws.Range("K34").Value = 1928.5834
$ws.Range("L34").Value = 25644456
$ws.Range("M34").Value = -1726.5834
$ws.Range("N34").Value = -25644860
$ws.Range("H99").Value = 1239.6666
$ws.Range("I99").Value = 1109.5
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1109.5
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 388.5
$ws.Range("N99").Value = -4496
$ws.Range("H107").Value = 1355.1111
$ws.Range("I107").Value = 1370.8572
$ws.Range("K107").Value = 1370.8572
$ws.Range("M107").Value = 549.1428000000001
$ws.Range("H126").Value = 1239.6666
$ws.Range("I126").Value = 1109.5
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 3328.5
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -858.5
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 7354343.5
$ws.Range("I132").Value = 10001221
$ws.Range("J132").Value = 1905.1666
$ws.Range("K132").Value = 30003663
$ws.Range("L132").Value = 5715.4998
$ws.Range("M132").Value = -30001133
$ws.Range("N132").Value = -10775.4998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1033.6154
$ws.Range("I5").Value = 379
$ws.Range("J5").Value = 2081
$ws.Range("K5").Value = 1137
$ws.Range("L5").Value = 6243
$ws.Range("M5").Value = -1025
$ws.Range("N5").Value = -6467
$ws.Range("H87").Value = 14969.6875
$ws.Range("I87").Value = 8557
$ws.Range("J87").Value = 19957.334
$ws.Range("K87").Value = 25671
$ws.Range("L87").Value = 59872.00199999999
$ws.Range("M87").Value = -24423
$ws.Range("N87").Value = -62368.00199999999
$ws.Range("H90").Value = 14969.6875
$ws.Range("I90").Value = 8557
$ws.Range("J90").Value = 19957.334
$ws.Range("K90").Value = 77013
$ws.Range("L90").Value = 179616.006
$ws.Range("M90").Value = -70773
$ws.Range("N90").Value = -192096.006
$ws.Range("H93").Value = 8700
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 8700
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 26100
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -29844
$ws.Range("H135").Value = 1033.6154
$ws.Range("I135").Value = 379
$ws.Range("J135").Value = 2081
$ws.Range("K135").Value = 3411
$ws.Range("L135").Value = 18729
$ws.Range("M135").Value = -876
$ws.Range("N135").Value = -23799

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2993.757
$ws.Range("I132").Value = 2233.9
$ws.Range("J132").Value = 4893.4
$ws.Range("K132").Value = 6701.700000000001
$ws.Range("L132").Value = 14680.2
$ws.Range("M132").Value = -4171.700000000001
$ws.Range("N132").Value = -19740.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9512.5
$ws.Range("I40").Value = 9516.666999999999
$ws.Range("J40").Value = 9500
$ws.Range("K40").Value = 9516.666999999999
$ws.Range("L40").Value = 9500
$ws.Range("M40").Value = -9380.666999999999
$ws.Range("N40").Value = -9772
$ws.Range("H139").Value = 53212.5
$ws.Range("J139").Value = 53212.5
$ws.Range("L139").Value = 53212.5
$ws.Range("N139").Value = -63492.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1004
$ws.Range("J2").Value = 1004
$ws.Range("L2").Value = 1004
$ws.Range("N2").Value = -1228
